$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns where old "jumlah" column (D) used to be, shifting the
# old D..H columns to F..J. This reproduces the column-width groupings for
# A, C, F, G, H, I, J exactly (the engine's EntireColumn.Insert() copies the
# width of a neighboring column only when that is unambiguous).
$ws.Range("D1:E1").EntireColumn.Insert()

# Row 1 headers
$ws.Range("A1").Value = "jenis_bantuan"
$ws.Range("B1").Value = "tanggal"
$ws.Range("C1").Value = "sasaran"
$ws.Range("D1").Value = "harga_satuan"
$ws.Range("E1").Value = "jumlah_penerima"
$ws.Range("F1").Value = "jumlah_bantuan"

# Row 2 data
$ws.Range("A2").Value = "Beras 3 kg"
$ws.Range("B2").Value = 45570
$ws.Range("C2").Value = "Lebak"
$ws.Range("D2").Value = 45000
$ws.Range("E2").Value = 1000
$ws.Range("F2").Value = "25 karung"

$ws.Range("F2").Select()
